# Weekly update: insert a new "Femacal de La Calera - Espárragos" record
# at row 31 (pushing the existing rows 31-42 down to 32-43), and fill in
# the new row with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 31, shifting rows 31:42 -> 32:43
$ws.Rows.Item(31).Insert()

# Populate the newly inserted row 31 with the new weekly record
$ws.Range("A31").Value = 3
$ws.Range("B31").Value = 'Femacal de La Calera'
$ws.Range("C31").Value = 'Coquimbo'
$ws.Range("D31").Value = 44875
$ws.Range("E31").Value = 5
$ws.Range("F31").Value = 300000000
$ws.Range("G31").Value = 'Espárragos'
$ws.Range("H31").Value = 'Verde'
$ws.Range("I31").Value = 'Primera'
$ws.Range("J31").Value = 2230
$ws.Range("K31").Value = 1400
$ws.Range("L31").Value = 1500
$ws.Range("M31").Value = 1450
$ws.Range("N31").Value = '$/kilo'
$ws.Range("O31").Value = 'Provincia de Quillota'
$ws.Range("P31").Value = 1450
$ws.Range("Q31").Value = 1
$ws.Range("R31").Value = 'Hortaliza'
